# Locomotive info is saved and loaded, and saving and loading refactored slightly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "S 8" sheet: insert a new row 2 with a new backlog story about removing
#    a locomotive from the track, pushing the existing stories down.
# ---------------------------------------------------------------------------
$wsS8 = $wb.Worksheets.Item("S 8")
$wsS8.Rows.Item(2).Insert()
$wsS8.Range("A2").Value = "I want a way to remove a locomotive from the track."
$wsS8.Rows.Item(2).RowHeight = 39

# ---------------------------------------------------------------------------
# 2) "Backlog" sheet: mark the four rows belonging to the "locomotive
#    save/load" story (rows 186-189) as Done, by putting an X in column C.
# ---------------------------------------------------------------------------
$wsBacklog = $wb.Worksheets.Item("Backlog")
$wsBacklog.Range("C186").Value = "X"
$wsBacklog.Range("C187").Value = "X"
$wsBacklog.Range("C188").Value = "X"
$wsBacklog.Range("C189").Value = "X"

# Update the view so the split pane/scroll position and the current
# selection reflect the newly-completed rows.
$wsBacklog.Activate()
$aw = $excel.ActiveWindow
$aw.SplitRow = 165
$wsBacklog.Range("A164").Select()
$wsBacklog.Range("C190").Select()
